$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 148 ("Hortaliza, Vega Monumental
# Concepcion - Aji" sheet). Excel shifts the existing rows 148:225 down to
# 149:226, preserving all of their values/styles, and the sheet's used range
# grows from A1:R225 to A1:R226.
$ws.Rows("148:148").Insert()

# Populate the newly inserted row 148 with the new weekly price-report entry.
$ws.Cells.Item(148, 1).Value2 = 11
$ws.Cells.Item(148, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(148, 3).Value2 = "Bíobío"
$ws.Cells.Item(148, 4).Value2 = 45146
$ws.Cells.Item(148, 5).Value2 = 8
$ws.Cells.Item(148, 6).Value2 = 100112021
$ws.Cells.Item(148, 7).Value2 = "Ají"
$ws.Cells.Item(148, 8).Value2 = "Inferno"
$ws.Cells.Item(148, 9).Value2 = "Primera"
$ws.Cells.Item(148, 10).Value2 = 40
$ws.Cells.Item(148, 11).Value2 = 17000
$ws.Cells.Item(148, 12).Value2 = 18000
$ws.Cells.Item(148, 13).Value2 = 17500
$ws.Cells.Item(148, 14).Value2 = "`$/caja 10 kilos"
$ws.Cells.Item(148, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(148, 16).Value2 = 1750
$ws.Cells.Item(148, 17).Value2 = 10
$ws.Cells.Item(148, 18).Value2 = "Hortaliza"
